$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "67.047.13"
$ws.Range("E2").Value2 = "  +1.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.276.61"
$ws.Range("E3").Value2 = "  -2.03%  "

$ws.Range("E4").Value2 = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "576.64"
$ws.Range("E5").Value2 = "  -1.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "171.81"
$ws.Range("E6").Value2 = "  -7.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.00"

$ws.Range("E8").Value2 = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "3.270.79"
$ws.Range("E9").Value2 = "  -2.08%  "

$ws.Range("E10").Value2 = "  -5.39%  "

$ws.Range("E11").Value2 = "  -2.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "44.83"
$ws.Range("E12").Value2 = "  -4.57%  "

$ws.Range("E13").Value2 = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "687.84"
$ws.Range("E14").Value2 = "  +2.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.802.91"
$ws.Range("E15").Value2 = "  -1.88%  "

$ws.Range("E16").Value2 = "  -3.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "67.175.38"
$ws.Range("E17").Value2 = "  +1.13%  "

$ws.Range("E18").Value2 = "  +0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.273.59"
$ws.Range("E19").Value2 = "  -2.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.16"
$ws.Range("E20").Value2 = "  -4.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "10.62"
$ws.Range("E21").Value2 = "  -4.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.880"
$ws.Range("E22").Value2 = "  -1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "16.85"
$ws.Range("E23").Value2 = "  -4.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "5.22"
$ws.Range("E24").Value2 = "  +3.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "97.98"
$ws.Range("E25").Value2 = "  -3.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "3.82"
$ws.Range("E26").Value2 = "  -5.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.63"
$ws.Range("E27").Value2 = "  -5.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "33.23"
$ws.Range("E28").Value2 = "  +2.69%  "

$ws.Range("E29").Value2 = "  -4.71%  "

$ws.Range("E30").Value2 = "  -2.63%  "

$ws.Range("E31").Value2 = "  -3.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "581.91"
$ws.Range("E32").Value2 = "  -5.38%  "

$ws.Range("E33").Value2 = "  -2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "3.811.57"
$ws.Range("E34").Value2 = "  -1.50%  "

$ws.Range("E35").Value2 = "  -3.37%  "

$ws.Range("E36").Value2 = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "55.35"
$ws.Range("E37").Value2 = "  -1.66%  "

$ws.Range("E38").Value2 = "  -16.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.128"
$ws.Range("E39").Value2 = "  -0.40%  "

$ws.Range("E40").Value2 = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "2.55"
$ws.Range("E41").Value2 = "  -4.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "31.26"
$ws.Range("E42").Value2 = "  -4.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0₃0653"
$ws.Range("E43").Value2 = "  -6.99%  "

$ws.Range("B44").Value2 = "TheGraph"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.323"
$ws.Range("E44").Value2 = "  -4.15%  "

$ws.Range("B45").Value2 = "Stacks"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.96"
$ws.Range("E45").Value2 = "  -7.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0401"
$ws.Range("E46").Value2 = "  -4.07%  "

$ws.Range("E47").Value2 = "  -0.02%  "

$ws.Range("E48").Value2 = "  -1.48%  "

$ws.Range("E49").Value2 = "  -1.11%  "

$ws.Range("E50").Value2 = "  +2.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "128.49"
$ws.Range("E51").Value2 = "  -0.40%  "
